$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = -12.452003414729768
$ws.Range("C2").Value = 6.1063278160393972
$ws.Range("D2").Value = 1.7946541736515533
$ws.Range("E2").Value = [double]"-3.4537018995118768E-2"

$ws.Range("B3").Value = 7.226494459758567
$ws.Range("C3").Value = 8.7873681415789751
$ws.Range("D3").Value = 15.667814657011764
$ws.Range("E3").Value = -7.3683557482878959

$ws.Range("B1:E3").Select()
